# Apply updated dSF (column F) values to Sheet1, per the repull/push/mean
# recalculation described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of row number -> new value for column F ("dSF")
$updates = @{
    2  = 0
    4  = 0
    8  = 0
    16 = -2
    22 = -2
    25 = 3
    27 = 4
    31 = 0
    34 = -2
    35 = -6
    36 = -4
    37 = 2
    38 = 1
    40 = 1
    41 = 3
    43 = 3
    45 = 3
    47 = -3
    48 = -2
    49 = 3
    54 = -6
    55 = -5
    57 = -7
    63 = -7
    67 = 0
    68 = -3
    70 = -6
    72 = -3
    74 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
